$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.46978126609108
$ws.Range("C2").Value = 7.902901218768448
$ws.Range("D2").Value = 3.808319522319407
$ws.Range("F2").Value = 20.79039870847368
$ws.Range("G2").Value = 3.59741220045176
$ws.Range("M2").Value = 19.96745915319175
$ws.Range("N2").Value = 16.48273794919262
$ws.Range("O2").Value = 18.38178476286791
$ws.Range("B3").Value = 9.965702345765218
$ws.Range("C3").Value = 7.642500481011892
$ws.Range("D3").Value = 3.761956772302128
$ws.Range("F3").Value = 20.69967622383982
$ws.Range("G3").Value = 3.59944262452683
$ws.Range("M3").Value = 19.36044299238007
$ws.Range("N3").Value = 16.55044511306207
$ws.Range("O3").Value = 18.38120398236523
$ws.Range("B4").Value = 9.643775446309027
$ws.Range("C4").Value = 7.476695600016303
$ws.Range("D4").Value = 3.732798885231843
$ws.Range("F4").Value = 20.65081845203585
$ws.Range("G4").Value = 3.600755320017431
$ws.Range("M4").Value = 18.9860236845136
$ws.Range("N4").Value = 16.59386810813637
$ws.Range("O4").Value = 18.38626834874109
$ws.Range("B5").Value = 9.509633033473424
$ws.Range("C5").Value = 7.407708067045711
$ws.Range("D5").Value = 3.720749647683164
$ws.Range("F5").Value = 20.63264588451259
$ws.Range("G5").Value = 3.601306907749528
$ws.Range("M5").Value = 18.83327292684265
$ws.Range("N5").Value = 16.61203060546836
$ws.Range("O5").Value = 18.38969332250826
$ws.Range("B6").Value = 9.487185792757261
$ws.Range("C6").Value = 7.39616903839074
$ws.Range("D6").Value = 3.718738980408371
$ws.Range("F6").Value = 20.62973368386143
$ws.Range("G6").Value = 3.601399505948613
$ws.Range("M6").Value = 18.8079061904981
$ws.Range("N6").Value = 16.6150747535388
$ws.Range("O6").Value = 18.39034417372639
$ws.Range("B7").Value = 9.64197807542662
$ws.Range("C7").Value = 7.47577086829547
$ws.Range("D7").Value = 3.732637053031752
$ws.Range("F7").Value = 20.6505663171143
$ws.Range("G7").Value = 3.600762691418983
$ws.Range("M7").Value = 18.98396397996257
$ws.Range("N7").Value = 16.59411115944174
$ws.Range("O7").Value = 18.38630903104105
$ws.Range("B8").Value = 10.29864587254525
$ws.Range("C8").Value = 7.814381572466564
$ws.Range("D8").Value = 3.792480956829943
$ws.Range("F8").Value = 20.7577070234889
$ws.Range("G8").Value = 3.598098624533174
$ws.Range("M8").Value = 19.75866515501102
$ws.Range("N8").Value = 16.50570047993233
$ws.Range("O8").Value = 18.38045851565702
$ws.Range("B9").Value = 11.48166594457186
$ws.Range("C9").Value = 8.428825341307199
$ws.Range("D9").Value = 3.904072373316259
$ws.Range("F9").Value = 21.02130231446118
$ws.Range("G9").Value = 3.593395589799213
$ws.Range("M9").Value = 21.25369898185258
$ws.Range("N9").Value = 16.34692059991866
$ws.Range("O9").Value = 18.41204173644015
$ws.Range("B10").Value = 12.28027604249875
$ws.Range("C10").Value = 8.846985739924818
$ws.Range("D10").Value = 3.982185661979862
$ws.Range("F10").Value = 21.24628231776407
$ws.Range("G10").Value = 3.590254440813121
$ws.Range("M10").Value = 22.32435082029166
$ws.Range("N10").Value = 16.23903396574347
$ws.Range("O10").Value = 18.46148901138047
$ws.Range("B11").Value = 12.62721029774221
$ws.Range("C11").Value = 9.029465891413162
$ws.Range("D11").Value = 4.016808205583964
$ws.Range("F11").Value = 21.35510360616356
$ws.Range("G11").Value = 3.588892910090961
$ws.Range("M11").Value = 22.8029901963447
$ws.Range("N11").Value = 16.19183019935221
$ws.Range("O11").Value = 18.48965894219354
$ws.Range("B12").Value = 12.75616480356128
$ws.Range("C12").Value = 9.097415861255904
$ws.Range("D12").Value = 4.029781977694014
$ws.Range("F12").Value = 21.39721176216435
$ws.Range("G12").Value = 3.588386966994227
$ws.Range("M12").Value = 22.9828539610874
$ws.Range("N12").Value = 16.17422280465591
$ws.Range("O12").Value = 18.50113856003073
$ws.Range("B13").Value = 12.72850084154657
$ws.Range("C13").Value = 9.082833376496058
$ws.Range("D13").Value = 4.026994030537407
$ws.Range("F13").Value = 21.38810351614328
$ws.Range("G13").Value = 3.58849550302044
$ws.Range("M13").Value = 22.94418139251418
$ws.Range("N13").Value = 16.17800299853198
$ws.Range("O13").Value = 18.49863016720889
$ws.Range("B14").Value = 12.63786835818583
$ws.Range("C14").Value = 9.035079450147785
$ws.Range("D14").Value = 4.017878345214707
$ws.Range("F14").Value = 21.358550021486
$ws.Range("G14").Value = 3.588851092964711
$ws.Range("M14").Value = 22.81781644109795
$ws.Range("N14").Value = 16.19037627517475
$ws.Range("O14").Value = 18.49058712453704
$ws.Range("B15").Value = 12.58203607452929
$ws.Range("C15").Value = 9.005677815408596
$ws.Range("D15").Value = 4.012276706398908
$ws.Range("F15").Value = 21.34056387802925
$ws.Range("G15").Value = 3.589070155818908
$ws.Range("M15").Value = 22.74022871168589
$ws.Range("N15").Value = 16.19799006385583
$ws.Range("O15").Value = 18.48576617518514
$ws.Range("B16").Value = 12.25726766465911
$ws.Range("C16").Value = 8.834901076154033
$ws.Range("D16").Value = 3.979904139988032
$ws.Range("F16").Value = 21.239298277144
$ws.Range("G16").Value = 3.590344771721219
$ws.Range("M16").Value = 22.29288585678573
$ws.Range("N16").Value = 16.24215639913101
$ws.Range("O16").Value = 18.45976194637218
$ws.Range("B17").Value = 12.05378949036235
$ws.Range("C17").Value = 8.728123408885546
$ws.Range("D17").Value = 3.959806833044169
$ws.Range("F17").Value = 21.17881223698146
$ws.Range("G17").Value = 3.591143931281045
$ws.Range("M17").Value = 22.01617039367353
$ws.Range("N17").Value = 16.26972975026814
$ws.Range("O17").Value = 18.44526053021351
$ws.Range("B18").Value = 11.93521787849956
$ws.Range("C18").Value = 8.665980842274516
$ws.Range("D18").Value = 3.948161932144596
$ws.Range("F18").Value = 21.14463390859841
$ws.Range("G18").Value = 3.591609933236149
$ws.Range("M18").Value = 21.85622929503941
$ws.Range("N18").Value = 16.28576575337618
$ws.Range("O18").Value = 18.43745426882442
$ws.Range("B19").Value = 11.89481003393513
$ws.Range("C19").Value = 8.644816827727283
$ws.Range("D19").Value = 3.944204675830856
$ws.Range("F19").Value = 21.1331676888136
$ws.Range("G19").Value = 3.591768805101647
$ws.Range("M19").Value = 21.80194744607115
$ws.Range("N19").Value = 16.29122564709775
$ws.Range("O19").Value = 18.43490311764338
$ws.Range("B20").Value = 12.07560964084605
$ws.Range("C20").Value = 8.739565592192839
$ws.Range("D20").Value = 3.961955114462758
$ws.Range("F20").Value = 21.18518799801525
$ws.Range("G20").Value = 3.591058202920418
$ws.Range("M20").Value = 22.04570947367225
$ws.Range("N20").Value = 16.26677626236312
$ws.Range("O20").Value = 18.44674892985486
$ws.Range("B21").Value = 12.6645555419021
$ws.Range("C21").Value = 9.049137460115638
$ws.Range("D21").Value = 4.020559608247987
$ws.Range("F21").Value = 21.36720643687151
$ws.Range("G21").Value = 3.588746386391792
$ws.Range("M21").Value = 22.85497185100868
$ws.Range("N21").Value = 16.18673469508073
$ws.Range("O21").Value = 18.49292755307027
$ws.Range("B22").Value = 13.03532463296727
$ws.Range("C22").Value = 9.244735590847302
$ws.Range("D22").Value = 4.0580597049379
$ws.Range("F22").Value = 21.49139534344579
$ws.Range("G22").Value = 3.587291640361147
$ws.Range("M22").Value = 23.37571951120724
$ws.Range("N22").Value = 16.13598207931531
$ws.Range("O22").Value = 18.52783991686195
$ws.Range("B23").Value = 12.83875204469976
$ws.Range("C23").Value = 9.140967824392957
$ws.Range("D23").Value = 4.038120422711809
$ws.Range("F23").Value = 21.42464579685837
$ws.Range("G23").Value = 3.588062944217397
$ws.Range("M23").Value = 23.09858678377433
$ws.Range("N23").Value = 16.16292766520012
$ws.Range("O23").Value = 18.50877514286555
$ws.Range("B24").Value = 12.06574969570294
$ws.Range("C24").Value = 8.734394928805303
$ws.Range("D24").Value = 3.960984158463162
$ws.Range("F24").Value = 21.18230365798975
$ws.Range("G24").Value = 3.591096940292571
$ws.Range("M24").Value = 22.03235750097607
$ws.Range("N24").Value = 16.26811096186683
$ws.Range("O24").Value = 18.44607437078061
$ws.Range("B25").Value = 11.17360653386256
$ws.Range("C25").Value = 8.268247518328108
$ws.Range("D25").Value = 3.874538778880024
$ws.Range("F25").Value = 20.94438945070537
$ws.Range("G25").Value = 3.594612457390045
$ws.Range("M25").Value = 20.85325030138399
$ws.Range("N25").Value = 16.38832564785153
$ws.Range("O25").Value = 18.39888445284592
